# Updated cryptos list values (price + volume) to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.927.96"
Set-TextValue $ws.Range("E2") "  +0.08%  "
Set-TextValue $ws.Range("D3") "1.894.37"
Set-TextValue $ws.Range("E3") "  -0.16%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "0.7749"
Set-TextValue $ws.Range("E5") "  -2.40%  "
Set-TextValue $ws.Range("D6") "244.72"
Set-TextValue $ws.Range("E6") "  +0.19%  "
Set-TextValue $ws.Range("D7") "1.001"
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "0.3135"
Set-TextValue $ws.Range("E8") "  -1.13%  "
Set-TextValue $ws.Range("D9") "25.82"
Set-TextValue $ws.Range("E9") "  +1.40%  "
Set-TextValue $ws.Range("D10") "0.07251"
Set-TextValue $ws.Range("E10") "  +0.93%  "
Set-TextValue $ws.Range("D11") "0.08993"
Set-TextValue $ws.Range("E11") "  +10.83%  "
Set-TextValue $ws.Range("D12") "0.7728"
Set-TextValue $ws.Range("E12") "  +0.42%  "
Set-TextValue $ws.Range("D13") "1.917.20"
Set-TextValue $ws.Range("E13") "  -0.26%  "
Set-TextValue $ws.Range("D14") "5.448"
Set-TextValue $ws.Range("E14") "  -3.34%  "
Set-TextValue $ws.Range("D15") "94.68"
Set-TextValue $ws.Range("E15") "  +2.20%  "
Set-TextValue $ws.Range("D16") "6.209"
Set-TextValue $ws.Range("E16") "  +0.43%  "
Set-TextValue $ws.Range("D17") "29.925.88"
Set-TextValue $ws.Range("E17") "  +0.00%  "
Set-TextValue $ws.Range("D18") "13.99"
Set-TextValue $ws.Range("E18") "  -0.05%  "
Set-TextValue $ws.Range("D19") "246.58"
Set-TextValue $ws.Range("E19") "  +0.57%  "
Set-TextValue $ws.Range("D20") "0.000007892"
Set-TextValue $ws.Range("E20") "  +1.43%  "
Set-TextValue $ws.Range("D21") "2.164.32"
Set-TextValue $ws.Range("E21") "  -0.02%  "
Set-TextValue $ws.Range("D22") "8.147"
Set-TextValue $ws.Range("E22") "  -2.19%  "
Set-TextValue $ws.Range("D23") "1.001"
Set-TextValue $ws.Range("E23") "  +0.01%  "
Set-TextValue $ws.Range("D24") "1.001"
Set-TextValue $ws.Range("E24") "  +0.04%  "
Set-TextValue $ws.Range("D25") "0.1588"
Set-TextValue $ws.Range("E25") "  -5.14%  "
Set-TextValue $ws.Range("D26") "9.546"
Set-TextValue $ws.Range("E26") "  +0.16%  "
Set-TextValue $ws.Range("D27") "162.56"
Set-TextValue $ws.Range("E27") "  -0.70%  "
Set-TextValue $ws.Range("E28") "  +0.34%  "
Set-TextValue $ws.Range("D29") "2.045"
Set-TextValue $ws.Range("E29") "  -1.34%  "
Set-TextValue $ws.Range("D30") "1.427"
Set-TextValue $ws.Range("E30") "  +1.62%  "
Set-TextValue $ws.Range("E31") "  +0.22%  "
Set-TextValue $ws.Range("D32") "4.524"
Set-TextValue $ws.Range("E32") "  +0.56%  "
Set-TextValue $ws.Range("D33") "4.117"
Set-TextValue $ws.Range("E33") "  +0.30%  "
Set-TextValue $ws.Range("D34") "0.05510"
Set-TextValue $ws.Range("E34") "  -2.32%  "
Set-TextValue $ws.Range("D35") "1.247"
Set-TextValue $ws.Range("E35") "  -2.67%  "
Set-TextValue $ws.Range("D36") "0.7543"
Set-TextValue $ws.Range("E36") "  +1.31%  "
Set-TextValue $ws.Range("D37") "1.001"
Set-TextValue $ws.Range("E37") "  -0.14%  "
Set-TextValue $ws.Range("D38") "2.707"
Set-TextValue $ws.Range("E38") "  +3.19%  "
Set-TextValue $ws.Range("D39") "0.01976"
Set-TextValue $ws.Range("E39") "  +2.02%  "
Set-TextValue $ws.Range("D40") "2.794"
Set-TextValue $ws.Range("E40") "  +0.23%  "
Set-TextValue $ws.Range("D41") "0.4510"
Set-TextValue $ws.Range("E41") "  +1.75%  "
Set-TextValue $ws.Range("D42") "74.20"
Set-TextValue $ws.Range("E42") "  -1.27%  "
Set-TextValue $ws.Range("D43") "6.085"
Set-TextValue $ws.Range("E43") "  +2.01%  "
Set-TextValue $ws.Range("D44") "1.092.43"
Set-TextValue $ws.Range("E44") "  -6.49%  "
Set-TextValue $ws.Range("D45") "0.8556"
Set-TextValue $ws.Range("E45") "  +0.16%  "
Set-TextValue $ws.Range("D46") "1.001"
Set-TextValue $ws.Range("E46") "  +0.05%  "
Set-TextValue $ws.Range("D47") "1.893"
Set-TextValue $ws.Range("E47") "  +0.15%  "
Set-TextValue $ws.Range("D48") "102.78"
Set-TextValue $ws.Range("E48") "  -1.90%  "
Set-TextValue $ws.Range("D49") "7.642"
Set-TextValue $ws.Range("E49") "  +2.06%  "
Set-TextValue $ws.Range("D50") "9.855"
Set-TextValue $ws.Range("E50") "  -2.10%  "
Set-TextValue $ws.Range("D51") "2.999"
Set-TextValue $ws.Range("E51") "  -0.12%  "
